$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("Nov 2" entry) is no longer the last row, so it switches from the
# "last row" date-only format to the regular date+time format used by all
# the other historical rows.
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (r15): Day/Chase/Bryce/Zach for 2025-11-03.
$ws.Range("A15").Value = 45964
$ws.Range("A15").NumberFormat = "YYYY-MM-DD"
$ws.Range("B15").Value = 29
$ws.Range("C15").Value = 38
$ws.Range("D15").Value = 36
